# Generate Report for Handback
# The 77c44467-84b7-4793-9c3c-6df43c77a886.md file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# per-locale sheets' Status + Latest Handback DateTime, and roll that up
# into the Overview sheet.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-03-20 14:43:55"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-03-20 14:44:01"
